$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.086.98'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.787.32'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.86%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.25'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.551'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.52%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.79'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.16%  '
$ws.Range("E9").Value = '  -2.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0710'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0932'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.045.14'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.795.79'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.89'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("E15").Value = '  -3.97%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.065.48'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.17'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.85'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.94%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.36'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0791'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.83'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.88%  '
$ws.Range("E23").Value = '  -4.70%  '
$ws.Range("E24").Value = '  -2.86%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.73'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.34'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.35%  '
$ws.Range("E27").Value = '  -2.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.113'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("E30").Value = '  -4.19%  '
$ws.Range("E31").Value = '  +0.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.66'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.52'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.89%  '
$ws.Range("E34").Value = '  -5.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.396.48'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.24%  '
$ws.Range("E36").Value = '  -0.86%  '
$ws.Range("E37").Value = '  -1.52%  '
$ws.Range("E38").Value = '  -3.43%  '
$ws.Range("E39").Value = '  +2.59%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("E41").Value = '  -5.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.70'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '77.88'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.56%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₆0145'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +14.90%  '
$ws.Range("E45").Value = '  +2.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.64'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +4.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0498'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("E48").Value = '  +1.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.944.43'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.89%  '
$ws.Range("E51").Value = '  +0.06%  '
